$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 'ECHO VENTURES GROUP LIMITED'
$ws.Range("B11").Value = 16455744
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '2025-05-18'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = 'active'
$ws.Range("E11").Value = 'Keyword'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '2025-05-18'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = '2025-05-18 21:42:48'

$ws.Range("A12").Value = 'ESLB INVESTMENTS LIMITED'
$ws.Range("B12").Value = 16455669
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '2025-05-18'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 'active'
$ws.Range("E12").Value = 'Keyword'
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = '2025-05-18'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = '2025-05-18 21:42:48'

$ws.Range("A13").Value = 'JENKINS VENTURES LTD'
$ws.Range("B13").Value = 16455788
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '2025-05-18'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = 'active'
$ws.Range("E13").Value = 'Keyword'
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = '2025-05-18'
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = '2025-05-18 21:42:49'

$ws.Range("A14").Value = 'MARIOS PROPERTY INVESTMENTS LTD'
$ws.Range("B14").Value = 16455816
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '2025-05-18'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = 'active'
$ws.Range("E14").Value = 'Keyword'
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = '2025-05-18'
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = '2025-05-18 21:42:50'

$ws.Range("A15").Value = 'LENDING CONSULTANCY LTD'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '16455471'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '2025-05-18'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = 'active'
$ws.Range("E15").Value = 'SIC'
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = '2025-05-18'
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = '2025-05-18 23:26:26'

$ws.Range("A16").Value = 'ECHO VENTURES GROUP LIMITED'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '16455744'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '2025-05-18'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 'active'
$ws.Range("E16").Value = 'Keyword'
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = '2025-05-18'
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = '2025-05-18 23:26:27'

$ws.Range("A17").Value = 'ESLB INVESTMENTS LIMITED'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '16455669'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = '2025-05-18'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = 'active'
$ws.Range("E17").Value = 'Keyword'
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = '2025-05-18'
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = '2025-05-18 23:26:27'

$ws.Range("A18").Value = 'JISA VENTURES LTD'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '16455405'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '2025-05-18'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = 'active'
$ws.Range("E18").Value = 'Keyword'
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = '2025-05-18'
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = '2025-05-18 23:26:27'

$ws.Range("A19").Value = 'TALLY M E VENTURES LIMITED'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '16455468'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '2025-05-18'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = 'active'
$ws.Range("E19").Value = 'Keyword'
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = '2025-05-18'
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = '2025-05-18 23:26:27'

$ws.Range("A20").Value = 'PERFICIENT VENTURES LTD'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '16455594'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '2025-05-18'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = 'active'
$ws.Range("E20").Value = 'Keyword'
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = '2025-05-18'
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = '2025-05-18 23:26:28'

$ws.Range("A21").Value = 'BLUEBOW TECHNOLOGIES LTD'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '16455597'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = '2025-05-18'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = 'active'
$ws.Range("E21").Value = 'SIC'
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = '2025-05-18'
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = '2025-05-18 23:26:28'

$ws.Range("A22").Value = 'JENKINS VENTURES LTD'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '16455788'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '2025-05-18'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = 'active'
$ws.Range("E22").Value = 'Keyword'
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = '2025-05-18'
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Value = '2025-05-18 23:26:28'

$ws.Range("A23").Value = 'BIEN DEVELOPMENTS LTD'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '16455494'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '2025-05-18'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = 'active'
$ws.Range("E23").Value = 'SIC'
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = '2025-05-18'
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value = '2025-05-18 23:26:28'

$ws.Range("A24").Value = 'ALPHA HAULAGE SOLUTIONS LTD'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '16455573'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = '2025-05-18'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = 'active'
$ws.Range("E24").Value = 'Keyword'
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = '2025-05-18'
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value = '2025-05-18 23:26:29'

$ws.Range("A25").Value = 'MARIOS PROPERTY INVESTMENTS LTD'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '16455816'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = '2025-05-18'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = 'active'
$ws.Range("E25").Value = 'Keyword'
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = '2025-05-18'
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = '2025-05-18 23:26:29'

$ws.Range("A26").Value = 'MARKOVIAN INVESTMENTS LIMITED'
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '16455443'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = '2025-05-18'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = 'active'
$ws.Range("E26").Value = 'Keyword'
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = '2025-05-18'
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = '2025-05-18 23:26:29'

$ws.Range("A27").Value = 'PARTNERS AMERICAN WHISKEY LTD'
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '16455528'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = '2025-05-18'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = 'active'
$ws.Range("E27").Value = 'Keyword'
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = '2025-05-18'
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value = '2025-05-18 23:26:29'
